# Insert a new weekly price record as row 32 on the active sheet.
# This pushes the existing rows 32..97 down to 33..98 (dimension grows
# from A1:R97 to A1:R98), matching the "Fruta / hortaliza, semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(32).Insert()

$ws.Range("A32").Value2 = 3
$ws.Range("B32").Value2 = "Femacal de La Calera"
$ws.Range("C32").Value2 = "Coquimbo"
$ws.Range("D32").Value2 = 44519
$ws.Range("E32").Value2 = 5
$ws.Range("F32").Value2 = 100112026
$ws.Range("G32").Value2 = "Haba"
$ws.Range("H32").Value2 = "Sin especificar"
$ws.Range("I32").Value2 = "Primera"
$ws.Range("J32").Value2 = 110
$ws.Range("K32").Value2 = 7000
$ws.Range("L32").Value2 = 7500
$ws.Range("M32").Value2 = 7273
$ws.Range("N32").Value2 = "`$/malla 25 kilos"
$ws.Range("O32").Value2 = "Provincia de Limarí"
$ws.Range("P32").Value2 = 291
$ws.Range("Q32").Value2 = 25
$ws.Range("R32").Value2 = "Hortaliza"
